$d = $word.ActiveDocument

$d.Content.Find.Execute("534÷8=66, 6", $true, $false, $false, $false, $false, $true, 1, $false, "978÷5=195, 3", 2)
$d.Content.Find.Execute("457÷5=91, 2", $true, $false, $false, $false, $false, $true, 1, $false, "289÷9=32, 1", 2)
$d.Content.Find.Execute("367÷5=73, 2", $true, $false, $false, $false, $false, $true, 1, $false, "518÷7=74, 0", 2)
$d.Content.Find.Execute("934÷2=467, 0", $true, $false, $false, $false, $false, $true, 1, $false, "352÷8=44, 0", 2)
$d.Content.Find.Execute("316÷6=52, 4", $true, $false, $false, $false, $false, $true, 1, $false, "531÷7=75, 6", 2)
$d.Content.Find.Execute("203÷6=33, 5", $true, $false, $false, $false, $false, $true, 1, $false, "551÷7=78, 5", 2)
$d.Content.Find.Execute("222÷6=37, 0", $true, $false, $false, $false, $false, $true, 1, $false, "588÷5=117, 3", 2)
$d.Content.Find.Execute("664÷4=166, 0", $true, $false, $false, $false, $false, $true, 1, $false, "117÷7=16, 5", 2)
$d.Content.Find.Execute("574÷7=82, 0", $true, $false, $false, $false, $false, $true, 1, $false, "590÷6=98, 2", 2)
$d.Content.Find.Execute("980÷6=163, 2", $true, $false, $false, $false, $false, $true, 1, $false, "687÷4=171, 3", 2)
$d.Content.Find.Execute("626÷4=156, 2", $true, $false, $false, $false, $false, $true, 1, $false, "894÷2=447, 0", 2)
$d.Content.Find.Execute("410÷6=68, 2", $true, $false, $false, $false, $false, $true, 1, $false, "953÷4=238, 1", 2)
$d.Content.Find.Execute("341÷7=48, 5", $true, $false, $false, $false, $false, $true, 1, $false, "154÷3=51, 1", 2)
$d.Content.Find.Execute("196÷5=39, 1", $true, $false, $false, $false, $false, $true, 1, $false, "826÷6=137, 4", 2)
$d.Content.Find.Execute("453÷3=151, 0", $true, $false, $false, $false, $false, $true, 1, $false, "187÷2=93, 1", 2)
$d.Content.Find.Execute("616÷2=308, 0", $true, $false, $false, $false, $false, $true, 1, $false, "982÷8=122, 6", 2)
$d.Content.Find.Execute("519÷8=64, 7", $true, $false, $false, $false, $false, $true, 1, $false, "275÷7=39, 2", 2)
$d.Content.Find.Execute("743÷6=123, 5", $true, $false, $false, $false, $false, $true, 1, $false, "696÷9=77, 3", 2)
$d.Content.Find.Execute("645÷6=107, 3", $true, $false, $false, $false, $false, $true, 1, $false, "186÷6=31, 0", 2)
$d.Content.Find.Execute("279÷6=46, 3", $true, $false, $false, $false, $false, $true, 1, $false, "740÷2=370, 0", 2)
$d.Content.Find.Execute("449÷9=49, 8", $true, $false, $false, $false, $false, $true, 1, $false, "638÷4=159, 2", 2)
$d.Content.Find.Execute("945÷6=157, 3", $true, $false, $false, $false, $false, $true, 1, $false, "276÷8=34, 4", 2)
$d.Content.Find.Execute("840÷9=93, 3", $true, $false, $false, $false, $false, $true, 1, $false, "829÷3=276, 1", 2)
$d.Content.Find.Execute("727÷4=181, 3", $true, $false, $false, $false, $false, $true, 1, $false, "623÷2=311, 1", 2)
$d.Content.Find.Execute("183÷7=26, 1", $true, $false, $false, $false, $false, $true, 1, $false, "341÷7=48, 5", 2)
